# Contest 19 and 20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31 -> Contest 19
$ws.Range("E31").Value = 70
$ws.Range("H31").Value = 40
$ws.Range("K31").Value = 30
$ws.Range("N31").Value = 100
$ws.Range("Q31").Value = 60
$ws.Range("T31").Value = 50
$ws.Range("W31").Value = 80
$ws.Range("Z31").Value = 20
$ws.Range("AC31").Value = 0

# Row 32 -> Contest 20
$ws.Range("E32").Value = 0
$ws.Range("H32").Value = 80
$ws.Range("K32").Value = 50
$ws.Range("N32").Value = 70
$ws.Range("Q32").Value = 60
$ws.Range("T32").Value = 100
$ws.Range("W32").Value = 40
$ws.Range("Z32").Value = 20
$ws.Range("AC32").Value = 30

$excel.CalculateFull()
